# Regenerate the handback report: statuses move from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamps are refreshed,
# and the (now resolved) error details are cleared out. Also widen a couple
# of columns that held the longer status text / shrink the now-empty error
# column.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664

# --- "zh-cn" sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("K2").Value = "2016-08-12 16:44:09"
$wsZhCn.Range("K3").Value = "2016-08-12 16:44:09"

$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

# --- "de-de" sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("K2").Value = "2016-08-12 16:44:20"
$wsDeDe.Range("K3").Value = "2016-08-12 16:44:20"

$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
